$d = $word.ActiveDocument
$d.Content.Find.Execute("71+8=", $true, $false, $false, $false, $false, $true, 1, $false, "60-35=", 2) | Out-Null
$d.Content.Find.Execute("57-6=", $true, $false, $false, $false, $false, $true, 1, $false, "28-6=", 2) | Out-Null
$d.Content.Find.Execute("82+6=", $true, $false, $false, $false, $false, $true, 1, $false, "70-12=", 2) | Out-Null
$d.Content.Find.Execute("74+7=", $true, $false, $false, $false, $false, $true, 1, $false, "56-9=", 2) | Out-Null
$d.Content.Find.Execute("98-23=", $true, $false, $false, $false, $false, $true, 1, $false, "16-10=", 2) | Out-Null
$d.Content.Find.Execute("88-41=", $true, $false, $false, $false, $false, $true, 1, $false, "41-27=", 2) | Out-Null
$d.Content.Find.Execute("4+1=", $true, $false, $false, $false, $false, $true, 1, $false, "10+64=", 2) | Out-Null
$d.Content.Find.Execute("81-46=", $true, $false, $false, $false, $false, $true, 1, $false, "71-49=", 2) | Out-Null
$d.Content.Find.Execute("92-75=", $true, $false, $false, $false, $false, $true, 1, $false, "77+20=", 2) | Out-Null
$d.Content.Find.Execute("1+52=", $true, $false, $false, $false, $false, $true, 1, $false, "2+72=", 2) | Out-Null
$d.Content.Find.Execute("97+1=", $true, $false, $false, $false, $false, $true, 1, $false, "31+63=", 2) | Out-Null
$d.Content.Find.Execute("1+50=", $true, $false, $false, $false, $false, $true, 1, $false, "87-83=", 2) | Out-Null
$d.Content.Find.Execute("22+23=", $true, $false, $false, $false, $false, $true, 1, $false, "51-3=", 2) | Out-Null
$d.Content.Find.Execute("78+3=", $true, $false, $false, $false, $false, $true, 1, $false, "43-31=", 2) | Out-Null
$d.Content.Find.Execute("15+84=", $true, $false, $false, $false, $false, $true, 1, $false, "41-3=", 2) | Out-Null
$d.Content.Find.Execute("69-26=", $true, $false, $false, $false, $false, $true, 1, $false, "57+33=", 2) | Out-Null
$d.Content.Find.Execute("84-56=", $true, $false, $false, $false, $false, $true, 1, $false, "76-7=", 2) | Out-Null
$d.Content.Find.Execute("74+23=", $true, $false, $false, $false, $false, $true, 1, $false, "86-44=", 2) | Out-Null
$d.Content.Find.Execute("11+34=", $true, $false, $false, $false, $false, $true, 1, $false, "52-23=", 2) | Out-Null
$d.Content.Find.Execute("12+85=", $true, $false, $false, $false, $false, $true, 1, $false, "88-44=", 2) | Out-Null
$d.Content.Find.Execute("21+56=", $true, $false, $false, $false, $false, $true, 1, $false, "17-15=", 2) | Out-Null
$d.Content.Find.Execute("20+72=", $true, $false, $false, $false, $false, $true, 1, $false, "24+11=", 2) | Out-Null
$d.Content.Find.Execute("64+2=", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=", 2) | Out-Null
$d.Content.Find.Execute("16+73=", $true, $false, $false, $false, $false, $true, 1, $false, "39-31=", 2) | Out-Null
$d.Content.Find.Execute("87-29=", $true, $false, $false, $false, $false, $true, 1, $false, "80+15=", 2) | Out-Null
$d.Content.Find.Execute("79-72=", $true, $false, $false, $false, $false, $true, 1, $false, "48-36=", 2) | Out-Null
$d.Content.Find.Execute("46-24=", $true, $false, $false, $false, $false, $true, 1, $false, "63-0=", 2) | Out-Null
$d.Content.Find.Execute("58-0=", $true, $false, $false, $false, $false, $true, 1, $false, "83-23=", 2) | Out-Null
$d.Content.Find.Execute("93-30=", $true, $false, $false, $false, $false, $true, 1, $false, "84-25=", 2) | Out-Null
$d.Content.Find.Execute("18+16=", $true, $false, $false, $false, $false, $true, 1, $false, "12+7=", 2) | Out-Null
$d.Content.Find.Execute("15+79=", $true, $false, $false, $false, $false, $true, 1, $false, "15+52=", 2) | Out-Null
$d.Content.Find.Execute("72+25=", $true, $false, $false, $false, $false, $true, 1, $false, "54-27=", 2) | Out-Null
$d.Content.Find.Execute("3+62=", $true, $false, $false, $false, $false, $true, 1, $false, "80-29=", 2) | Out-Null
$d.Content.Find.Execute("63+33=", $true, $false, $false, $false, $false, $true, 1, $false, "47+51=", 2) | Out-Null
$d.Content.Find.Execute("30+36=", $true, $false, $false, $false, $false, $true, 1, $false, "4+80=", 2) | Out-Null
$d.Content.Find.Execute("89-58=", $true, $false, $false, $false, $false, $true, 1, $false, "0+42=", 2) | Out-Null
$d.Content.Find.Execute("65+16=", $true, $false, $false, $false, $false, $true, 1, $false, "18-10=", 2) | Out-Null
$d.Content.Find.Execute("76-68=", $true, $false, $false, $false, $false, $true, 1, $false, "75+5=", 2) | Out-Null
$d.Content.Find.Execute("23+51=", $true, $false, $false, $false, $false, $true, 1, $false, "52+13=", 2) | Out-Null
$d.Content.Find.Execute("26+6=", $true, $false, $false, $false, $false, $true, 1, $false, "96-60=", 2) | Out-Null
$d.Content.Find.Execute("94-57=", $true, $false, $false, $false, $false, $true, 1, $false, "63-59=", 2) | Out-Null
$d.Content.Find.Execute("94-43=", $true, $false, $false, $false, $false, $true, 1, $false, "31+64=", 2) | Out-Null
$d.Content.Find.Execute("28-20=", $true, $false, $false, $false, $false, $true, 1, $false, "86-67=", 2) | Out-Null
$d.Content.Find.Execute("35+1=", $true, $false, $false, $false, $false, $true, 1, $false, "60-31=", 2) | Out-Null
$d.Content.Find.Execute("50+32=", $true, $false, $false, $false, $false, $true, 1, $false, "0+69=", 2) | Out-Null
$d.Content.Find.Execute("3+71=", $true, $false, $false, $false, $false, $true, 1, $false, "47-20=", 2) | Out-Null
$d.Content.Find.Execute("42-42=", $true, $false, $false, $false, $false, $true, 1, $false, "95-73=", 2) | Out-Null
$d.Content.Find.Execute("76+7=", $true, $false, $false, $false, $false, $true, 1, $false, "45+7=", 2) | Out-Null
$d.Content.Find.Execute("29+68=", $true, $false, $false, $false, $false, $true, 1, $false, "78+20=", 2) | Out-Null
$d.Content.Find.Execute("58-1=", $true, $false, $false, $false, $false, $true, 1, $false, "51-44=", 2) | Out-Null
$d.Content.Find.Execute("33+65=", $true, $false, $false, $false, $false, $true, 1, $false, "99-45=", 2) | Out-Null
$d.Content.Find.Execute("37-26=", $true, $false, $false, $false, $false, $true, 1, $false, "51-1=", 2) | Out-Null
$d.Content.Find.Execute("81-58=", $true, $false, $false, $false, $false, $true, 1, $false, "31+25=", 2) | Out-Null
$d.Content.Find.Execute("20+55=", $true, $false, $false, $false, $false, $true, 1, $false, "33+37=", 2) | Out-Null
$d.Content.Find.Execute("85-83=", $true, $false, $false, $false, $false, $true, 1, $false, "36+36=", 2) | Out-Null
$d.Content.Find.Execute("83-22=", $true, $false, $false, $false, $false, $true, 1, $false, "18+45=", 2) | Out-Null
$d.Content.Find.Execute("23-11=", $true, $false, $false, $false, $false, $true, 1, $false, "53-19=", 2) | Out-Null
$d.Content.Find.Execute("83+1=", $true, $false, $false, $false, $false, $true, 1, $false, "86+13=", 2) | Out-Null
$d.Content.Find.Execute("54+20=", $true, $false, $false, $false, $false, $true, 1, $false, "57+33=", 2) | Out-Null
$d.Content.Find.Execute("32+65=", $true, $false, $false, $false, $false, $true, 1, $false, "9+54=", 2) | Out-Null
$d.Content.Find.Execute("7+26=", $true, $false, $false, $false, $false, $true, 1, $false, "72+26=", 2) | Out-Null
$d.Content.Find.Execute("96-85=", $true, $false, $false, $false, $false, $true, 1, $false, "57-30=", 2) | Out-Null
$d.Content.Find.Execute("17+0=", $true, $false, $false, $false, $false, $true, 1, $false, "3+10=", 2) | Out-Null
$d.Content.Find.Execute("42+5=", $true, $false, $false, $false, $false, $true, 1, $false, "78-58=", 2) | Out-Null
$d.Content.Find.Execute("47+5=", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=", 2) | Out-Null
$d.Content.Find.Execute("90-27=", $true, $false, $false, $false, $false, $true, 1, $false, "99-67=", 2) | Out-Null
$d.Content.Find.Execute("81-75=", $true, $false, $false, $false, $false, $true, 1, $false, "63-61=", 2) | Out-Null
$d.Content.Find.Execute("55+33=", $true, $false, $false, $false, $false, $true, 1, $false, "85-30=", 2) | Out-Null
$d.Content.Find.Execute("48+12=", $true, $false, $false, $false, $false, $true, 1, $false, "87+11=", 2) | Out-Null
$d.Content.Find.Execute("79-56=", $true, $false, $false, $false, $false, $true, 1, $false, "30+32=", 2) | Out-Null
$d.Content.Find.Execute("97-57=", $true, $false, $false, $false, $false, $true, 1, $false, "56-32=", 2) | Out-Null
$d.Content.Find.Execute("36-6=", $true, $false, $false, $false, $false, $true, 1, $false, "17-1=", 2) | Out-Null
$d.Content.Find.Execute("66-46=", $true, $false, $false, $false, $false, $true, 1, $false, "52+27=", 2) | Out-Null
$d.Content.Find.Execute("69-23=", $true, $false, $false, $false, $false, $true, 1, $false, "84-61=", 2) | Out-Null
$d.Content.Find.Execute("79+12=", $true, $false, $false, $false, $false, $true, 1, $false, "82-23=", 2) | Out-Null
$d.Content.Find.Execute("12+70=", $true, $false, $false, $false, $false, $true, 1, $false, "83-37=", 2) | Out-Null
$d.Content.Find.Execute("94-9=", $true, $false, $false, $false, $false, $true, 1, $false, "14+33=", 2) | Out-Null
$d.Content.Find.Execute("73-39=", $true, $false, $false, $false, $false, $true, 1, $false, "91-49=", 2) | Out-Null
$d.Content.Find.Execute("98-65=", $true, $false, $false, $false, $false, $true, 1, $false, "47-13=", 2) | Out-Null
$d.Content.Find.Execute("32+11=", $true, $false, $false, $false, $false, $true, 1, $false, "82-64=", 2) | Out-Null
$d.Content.Find.Execute("17+37=", $true, $false, $false, $false, $false, $true, 1, $false, "10+6=", 2) | Out-Null
$d.Content.Find.Execute("73-56=", $true, $false, $false, $false, $false, $true, 1, $false, "26+4=", 2) | Out-Null
$d.Content.Find.Execute("94-76=", $true, $false, $false, $false, $false, $true, 1, $false, "63-44=", 2) | Out-Null
$d.Content.Find.Execute("77-22=", $true, $false, $false, $false, $false, $true, 1, $false, "8+62=", 2) | Out-Null
$d.Content.Find.Execute("6+11=", $true, $false, $false, $false, $false, $true, 1, $false, "15+61=", 2) | Out-Null
$d.Content.Find.Execute("67-3=", $true, $false, $false, $false, $false, $true, 1, $false, "42-7=", 2) | Out-Null
$d.Content.Find.Execute("23+28=", $true, $false, $false, $false, $false, $true, 1, $false, "51+43=", 2) | Out-Null
$d.Content.Find.Execute("60-26=", $true, $false, $false, $false, $false, $true, 1, $false, "94-25=", 2) | Out-Null
$d.Content.Find.Execute("18+29=", $true, $false, $false, $false, $false, $true, 1, $false, "43+41=", 2) | Out-Null
$d.Content.Find.Execute("82-81=", $true, $false, $false, $false, $false, $true, 1, $false, "50-1=", 2) | Out-Null
$d.Content.Find.Execute("97-52=", $true, $false, $false, $false, $false, $true, 1, $false, "48-20=", 2) | Out-Null
$d.Content.Find.Execute("62-52=", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=", 2) | Out-Null
$d.Content.Find.Execute("15-4=", $true, $false, $false, $false, $false, $true, 1, $false, "11+27=", 2) | Out-Null
$d.Content.Find.Execute("69+22=", $true, $false, $false, $false, $false, $true, 1, $false, "26-26=", 2) | Out-Null
$d.Content.Find.Execute("85-20=", $true, $false, $false, $false, $false, $true, 1, $false, "97-74=", 2) | Out-Null
$d.Content.Find.Execute("98-9=", $true, $false, $false, $false, $false, $true, 1, $false, "74+25=", 2) | Out-Null
$d.Content.Find.Execute("83-75=", $true, $false, $false, $false, $false, $true, 1, $false, "45-14=", 2) | Out-Null
$d.Content.Find.Execute("42+56=", $true, $false, $false, $false, $false, $true, 1, $false, "88-40=", 2) | Out-Null
$d.Content.Find.Execute("35-15=", $true, $false, $false, $false, $false, $true, 1, $false, "74+2=", 2) | Out-Null
$d.Content.Find.Execute("87-6=", $true, $false, $false, $false, $false, $true, 1, $false, "93-51=", 2) | Out-Null
